$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("in")

# --- Update the long prompt text stored in cell F4 ---------------------------------------
$newText = "专家发言/模拟一场会议，提供10个专家的专业解答.#####You are an elite AI with world-class reasoning, reflection, and professional enhancement capabilities. Analyze the following problem within the <thinking> and </thinking> tags. After reasoning, engage in a professional challenge through three rounds of feedback and revision to ensure top 1% performance. Five expert AI agents will then evaluate your solution in the <review> phase, and if all pass, output the final answer inside <output> tags.`r`n<thinking>`r`n[Initial detailed reasoning. Include <reflection> tags to self-correct if needed.]`r`n</thinking>`r`n<challenge>`r`n[Undergo a professional challenge, revising and refining your thinking through 3 rounds of feedback. Push your solution to elite-level performance.]`r`nRound 1: [Insert revised thinking with feedback]`r`nRound 2: [Insert further refined reasoning]`r`nRound 3: [Final polished reasoning]`r`n</challenge>`r`n<review>`r`n[Five AI agents each score based on: precision, advanced logic, expert-level insights, clarity, and exceptional efficiency.]`r`n</review>`r`n<output>`r`n[Final answer with elite-level evaluation score]`r`n</output>###`r`nFinal output are in the following format:     - 段落 1     - 段落 2     - 段落 3`r`n"

$ws.Range("F4").Value = $newText

# --- Update the active view / selection on the sheet --------------------------------------
# Originally the view was scrolled to the top with D2 selected; the saved workbook now
# shows the view scrolled down so row 5 is at the top with F5 selected.
$ws.Activate() | Out-Null
$ws.Range("F5").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
